$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.396.56"
$ws.Range("E2").Value = "  -0.89%  "
$ws.Range("D3").Value = "1.639.98"
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "'211.83"
$ws.Range("E5").Value = "  -1.55%  "
$ws.Range("D6").Value = "'0.528"
$ws.Range("E6").Value = "  +3.57%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("D8").Value = "'22.99"
$ws.Range("E8").Value = "  -4.36%  "
$ws.Range("D9").Value = "'0.257"
$ws.Range("E9").Value = "  -2.26%  "
$ws.Range("E10").Value = "  -1.93%  "
$ws.Range("D11").Value = "'0.0891"
$ws.Range("E11").Value = "  +1.38%  "
$ws.Range("D12").Value = "1.870.98"
$ws.Range("E12").Value = "  -1.75%  "
$ws.Range("D13").Value = "1.634.28"
$ws.Range("E13").Value = "  -2.02%  "
$ws.Range("B14").Value = "Polygon"
$ws.Range("C14").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D14").Value = "'0.568"
$ws.Range("E14").Value = "  +1.50%  "
$ws.Range("B15").Value = "Polkadot"
$ws.Range("C15").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D15").Value = "'4.04"
$ws.Range("E15").Value = "  -2.59%  "
$ws.Range("D16").Value = "'64.39"
$ws.Range("E16").Value = "  -3.60%  "
$ws.Range("D17").Value = "27.368.75"
$ws.Range("E17").Value = "  -0.89%  "
$ws.Range("D18").Value = "'229.27"
$ws.Range("E18").Value = "  -5.90%  "
$ws.Range("E19").Value = "  -1.45%  "
$ws.Range("D20").Value = "'7.58"
$ws.Range("E20").Value = "  -1.62%  "
$ws.Range("E21").Value = "  +0.08%  "
$ws.Range("D22").Value = "'4.33"
$ws.Range("E22").Value = "  -3.88%  "
$ws.Range("D23").Value = "'9.73"
$ws.Range("E23").Value = "  +4.18%  "
$ws.Range("E24").Value = "  -0.64%  "
$ws.Range("D25").Value = "'146.75"
$ws.Range("E25").Value = "  -0.41%  "
$ws.Range("D26").Value = "'6.99"
$ws.Range("E26").Value = "  -3.29%  "
$ws.Range("D27").Value = "'0.113"
$ws.Range("E27").Value = "  +1.06%  "
$ws.Range("E28").Value = "  -0.07%  "
$ws.Range("D29").Value = "'15.54"
$ws.Range("E29").Value = "  -6.00%  "
$ws.Range("E30").Value = "  -3.95%  "
$ws.Range("E31").Value = "  -3.87%  "
$ws.Range("D33").Value = "'3.15"
$ws.Range("E33").Value = "  +0.94%  "
$ws.Range("D34").Value = "1.411.51"
$ws.Range("E34").Value = "  -4.10%  "
$ws.Range("E35").Value = "  -0.02%  "
$ws.Range("E36").Value = "  -0.22%  "
$ws.Range("D37").Value = "'0.563"
$ws.Range("E37").Value = "  -2.35%  "
$ws.Range("D38").Value = "'0.882"
$ws.Range("E38").Value = "  -5.21%  "
$ws.Range("E39").Value = "  -3.64%  "
$ws.Range("D40").Value = "'1.02"
$ws.Range("E40").Value = "  +0.61%  "
$ws.Range("E41").Value = "  +0.02%  "
$ws.Range("D42").Value = "'2.46"
$ws.Range("E42").Value = "  -2.75%  "
$ws.Range("E43").Value = "  +1.48%  "
$ws.Range("E44").Value = "  +0.95%  "
$ws.Range("D45").Value = "'0.800"
$ws.Range("E45").Value = "  +1.62%  "
$ws.Range("D46").Value = "'64.41"
$ws.Range("E46").Value = "  -7.54%  "
$ws.Range("D47").Value = "1.780.72"
$ws.Range("E47").Value = "  -1.69%  "
$ws.Range("E48").Value = "  -5.29%  "
$ws.Range("D49").Value = "'87.78"
$ws.Range("E49").Value = "  -1.69%  "
$ws.Range("B51").Value = "Algorand"
$ws.Range("C51").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D51").Value = "'0.0988"
$ws.Range("E51").Value = "  -4.15%  "
